# RoboCUP_SSL_Hardware_BOM.xlsx edit
# "added an object detection sensor that we could use."
#
# Populates the IMU BOM sheet with the WSEN-ISDS 6-axis IMU part and
# renames/populates the "Active IR (find obstical) BOM" sheet into an
# "Obstical detection BOM" sheet carrying the OPT8241NBN sensor data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the obstacle-detection sheet.
# ---------------------------------------------------------------------
$wsObst = $wb.Worksheets.Item("Active IR (find obstical) BOM")
$wsObst.Name = "Obstical detection BOM"

$wsImu = $wb.Worksheets.Item("IMU BOM")
$wsComplete = $wb.Worksheets.Item("Complete BOM")
$wsBreak = $wb.Worksheets.Item("Break Beam BOM")

# ---------------------------------------------------------------------
# 2. IMU BOM sheet — WSEN-ISDS 6 Axis IMU
#    (order of the first-seen literal strings controls shared-string
#    index allocation, so write B2 before B1 etc. to mirror the
#    original authoring order.)
# ---------------------------------------------------------------------
$wsImu.Hyperlinks.Add($wsImu.Range("B2"), "https://www.we-online.com/en/components/products/WSEN-ISDS?sq=2536030320001", "2536030320001", "", "https://www.we-online.com/en/components/products/WSEN-ISDS?sq=2536030320001 - 2536030320001")
$wsImu.Range("B2").Value = "`u{2063}WSEN-ISDS 6 Axis IMU (Inertial Measurement Unit) & EV-Kits | Sensors | W`u{00FC}rth Elektronik Product Catalog (we-online.com)"
$wsImu.Range("B2").Style = "Hyperlink"

$wsImu.Range("B1").Value = "WSEN-ISDS 6 Axis IMU"
$wsImu.Range("B3").Value = "3,0*2,5*0,86"
$wsImu.Range("B4").Value = 1
$wsImu.Range("B5").Value = "Not applicable (sent free of charge by W`u{00FC}rth Electronics)"

# ---------------------------------------------------------------------
# 3. Relabel the "Data sheet" header to "Data sheet OR url" — used both
#    on the Complete BOM summary sheet and the Obstical detection sheet.
# ---------------------------------------------------------------------
$wsObst.Range("B1").Value = "OPT8241NBN"
$wsComplete.Range("A3").Value = "Data sheet OR url"
$wsObst.Range("A2").Value = "Data sheet OR url"

# ---------------------------------------------------------------------
# 4. Obstical detection BOM sheet — OPT8241NBN
# ---------------------------------------------------------------------
$wsObst.Hyperlinks.Add($wsObst.Range("B2"), "https://www.mouser.se/ProductDetail/Texas-Instruments/OPT8241NBN?qs=cGEy3R83DS%2FxFMUAL%252BoBvw%3D%3D", "", "", "https://www.mouser.se/ProductDetail/Texas-Instruments/OPT8241NBN?qs=cGEy3R83DS%2FxFMUAL%252BoBvw%3D%3D")
$wsObst.Range("B2").Value = "OPT8241NBN Texas Instruments | Mouser Sverige"
$wsObst.Range("B2").Style = "Hyperlink"

$wsObst.Range("B3").Value = "7,9*8,8*0,8"
$wsObst.Range("B4").Value = 1
$wsObst.Range("B5").Value = 605.13

# ---------------------------------------------------------------------
# 5. Break Beam BOM — drop the stray total-cost formula in B6 (no
#    longer applicable once the other BOM sheets are populated).
# ---------------------------------------------------------------------
$wsBreak.Range("B6").ClearContents()

# ---------------------------------------------------------------------
# 6. Column widths — best-fit widened once the long part names / URLs
#    were entered.
# ---------------------------------------------------------------------
$wsImu.Columns.Item(2).ColumnWidth = 102.49047851562500
$wsObst.Columns.Item(2).ColumnWidth = 3.0833333333333335
$wsComplete.Columns.Item(5).ColumnWidth = 102.49047851562500

# ---------------------------------------------------------------------
# 7. Sheet-view / selection bookkeeping to match the saved state.
# ---------------------------------------------------------------------
$wsBreak.Range("B6").Select()
$wsImu.Range("B24").Select()
$wsObst.Range("B6").Select()

$wsComplete.Activate()
$wsComplete.Range("D34").Select()

Write-Output "done"
